$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Vol Product Lines have been updated" -> EffectiveDate (F) and PreviousExpDate (I)
# dates move from 02/27/2023 to 07/15/2023 for both data rows.
$ws.Range("F2").Value = "07152023"
$ws.Range("I2").Value = "07152023"
$ws.Range("F3").Value = "07152023"
$ws.Range("I3").Value = "07152023"

# Update the saved view/selection state of the sheet.
$ws.Range("H5").Select()
